$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3526261448860168
$ws.Range("B1").Value = 3.745977878570557
$ws.Range("C1").Value = 5.874599933624268
$ws.Range("D1").Value = 1.641469120979309
$ws.Range("E1").Value = 0.9857849478721619
